{"js": "// Author's edit appended two new paragraphs at the end of the document\n// body (after the existing \"Some text in word!\" paragraph): one blank\n// paragraph, followed by a paragraph with the new sentence.\nconst body = context.document.body;\n\n// New blank paragraph at the end of the body.\nbody.insertParagraph(\"\", \"End\");\n\n// New paragraph with the added sentence, also appended at the end.\nbody.insertParagraph(\"Added some contetn to word file\", \"End\");\n\nawait context.sync();\n", "ps1": "# Add a blank paragraph and a new paragraph of text to the end of the\n# document body, mirroring the author's edit (two new paragraphs appended\n# after \"Some text in word!\").\n$d = $word.ActiveDocument\n\n# New blank paragraph right after the existing content.\n$r = $d.Paragraphs.Last.Range\n$r.InsertParagraphAfter()\n\n# New paragraph containing the added sentence.\n$r2 = $d.Paragraphs.Last.Range\n$r2.InsertParagraphAfter()\n$r3 = $d.Paragraphs.Last.Range\n$r3.Text = \"Added some contetn to word file\"\n"}
